$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1047397243"
$ws.Range("D16").Value = "MAIRA ALEJANDRA MARTINEZ CASTELLAR"
$ws.Range("E16").Value = "1607"
$ws.Range("F16").Value = 27578
$ws.Range("G16").Value = 689454
$ws.Range("C17").Value = "1044800187"
$ws.Range("D17").Value = "OSCAR IVAN IMITOLA HERRERA"
$ws.Range("E17").Value = "1608"
$ws.Range("F17").Value = 120000
$ws.Range("G17").Value = 3000000
$ws.Range("C18").Value = "1047397243"
$ws.Range("D18").Value = "MAIRA ALEJANDRA MARTINEZ CASTELLAR"
$ws.Range("E18").Value = "1608"
$ws.Range("F18").Value = 27578
$ws.Range("G18").Value = 689454
$ws.Range("C19").Value = "1050963899"
$ws.Range("D19").Value = "MARIA JOSE ROJAS HURTADO"
$ws.Range("E19").Value = "1608"
$ws.Range("F19").Value = 27578
$ws.Range("G19").Value = 689454
$ws.Range("C20").Value = "1044800187"
$ws.Range("D20").Value = "OSCAR IVAN IMITOLA HERRERA"
$ws.Range("E20").Value = "1609"
$ws.Range("F20").Value = 120000
$ws.Range("G20").Value = 3000000
$ws.Range("C21").Value = "1047397243"
$ws.Range("D21").Value = "MAIRA ALEJANDRA MARTINEZ CASTELLAR"
$ws.Range("E21").Value = "1609"
$ws.Range("F21").Value = 27578
$ws.Range("G21").Value = 689454
$ws.Range("C22").Value = "1050963899"
$ws.Range("D22").Value = "MARIA JOSE ROJAS HURTADO"
$ws.Range("E22").Value = "1609"
$ws.Range("F22").Value = 27578
$ws.Range("G22").Value = 689454
$ws.Range("C23").Value = "1044800187"
$ws.Range("D23").Value = "OSCAR IVAN IMITOLA HERRERA"
$ws.Range("E23").Value = "1610"
$ws.Range("F23").Value = 120000
$ws.Range("G23").Value = 3000000
$ws.Range("C24").Value = "1047397243"
$ws.Range("D24").Value = "MAIRA ALEJANDRA MARTINEZ CASTELLAR"
$ws.Range("E24").Value = "1610"
$ws.Range("F24").Value = 27578
$ws.Range("G24").Value = 689454
$ws.Range("C25").Value = "1050963899"
$ws.Range("D25").Value = "MARIA JOSE ROJAS HURTADO"
$ws.Range("E25").Value = "1610"
$ws.Range("F25").Value = 27578
$ws.Range("G25").Value = 689454
$ws.Range("C26").Value = "1044800187"
$ws.Range("D26").Value = "OSCAR IVAN IMITOLA HERRERA"
$ws.Range("E26").Value = "1611"
$ws.Range("F26").Value = 120000
$ws.Range("G26").Value = 3000000
$ws.Range("C27").Value = "1047397243"
$ws.Range("D27").Value = "MAIRA ALEJANDRA MARTINEZ CASTELLAR"
$ws.Range("E27").Value = "1611"
$ws.Range("F27").Value = 27578
$ws.Range("G27").Value = 689454
$ws.Range("C28").Value = "1050963899"
$ws.Range("D28").Value = "MARIA JOSE ROJAS HURTADO"
$ws.Range("E28").Value = "1611"
$ws.Range("F28").Value = 27578
$ws.Range("G28").Value = 689454
$ws.Range("C29").Value = "1044800187"
$ws.Range("D29").Value = "OSCAR IVAN IMITOLA HERRERA"
$ws.Range("E29").Value = "1612"
$ws.Range("F29").Value = 120000
$ws.Range("G29").Value = 3000000
$ws.Range("C30").Value = "1047397243"
$ws.Range("D30").Value = "MAIRA ALEJANDRA MARTINEZ CASTELLAR"
$ws.Range("E30").Value = "1612"
$ws.Range("F30").Value = 27578
$ws.Range("G30").Value = 689454
$ws.Range("C31").Value = "1050963899"
$ws.Range("D31").Value = "MARIA JOSE ROJAS HURTADO"
$ws.Range("E31").Value = "1612"
$ws.Range("F31").Value = 27578
$ws.Range("G31").Value = 689454
$ws.Range("C32").Value = "1044800187"
$ws.Range("D32").Value = "OSCAR IVAN IMITOLA HERRERA"
$ws.Range("E32").Value = "1701"
$ws.Range("F32").Value = 120000
$ws.Range("G32").Value = 3000000
$ws.Range("C33").Value = "1047397243"
$ws.Range("D33").Value = "MAIRA ALEJANDRA MARTINEZ CASTELLAR"
$ws.Range("E33").Value = "1701"
$ws.Range("F33").Value = 27578
$ws.Range("G33").Value = 689454
$ws.Range("C34").Value = "1050963899"
$ws.Range("D34").Value = "MARIA JOSE ROJAS HURTADO"
$ws.Range("E34").Value = "1701"
$ws.Range("F34").Value = 27578
$ws.Range("G34").Value = 689454
$ws.Range("C35").Value = "1044800187"
$ws.Range("D35").Value = "OSCAR IVAN IMITOLA HERRERA"
$ws.Range("E35").Value = "1702"
$ws.Range("F35").Value = 120000
$ws.Range("G35").Value = 3000000
$ws.Range("C36").Value = "1047397243"
$ws.Range("D36").Value = "MAIRA ALEJANDRA MARTINEZ CASTELLAR"
$ws.Range("E36").Value = "1702"
$ws.Range("F36").Value = 27578
$ws.Range("G36").Value = 689454
$ws.Range("C37").Value = "1050963899"
$ws.Range("D37").Value = "MARIA JOSE ROJAS HURTADO"
$ws.Range("E37").Value = "1702"
$ws.Range("F37").Value = 27578
$ws.Range("G37").Value = 689454
$ws.Range("C38").Value = "1044800187"
$ws.Range("D38").Value = "OSCAR IVAN IMITOLA HERRERA"
$ws.Range("E38").Value = "1703"
$ws.Range("F38").Value = 120000
$ws.Range("G38").Value = 3000000
$ws.Range("C39").Value = "1047397243"
$ws.Range("D39").Value = "MAIRA ALEJANDRA MARTINEZ CASTELLAR"
$ws.Range("E39").Value = "1703"
$ws.Range("F39").Value = 27578
$ws.Range("G39").Value = 689454
$ws.Range("C40").Value = "1050963899"
$ws.Range("D40").Value = "MARIA JOSE ROJAS HURTADO"
$ws.Range("E40").Value = "1703"
$ws.Range("F40").Value = 27578
$ws.Range("G40").Value = 689454
$ws.Range("C41").Value = "1044800187"
$ws.Range("D41").Value = "OSCAR IVAN IMITOLA HERRERA"
$ws.Range("E41").Value = "1704"
$ws.Range("F41").Value = 120000
$ws.Range("G41").Value = 3000000
$ws.Range("C42").Value = "1047397243"
$ws.Range("D42").Value = "MAIRA ALEJANDRA MARTINEZ CASTELLAR"
$ws.Range("E42").Value = "1704"
$ws.Range("F42").Value = 27578
$ws.Range("G42").Value = 689454
$ws.Range("C43").Value = "1050963899"
$ws.Range("D43").Value = "MARIA JOSE ROJAS HURTADO"
$ws.Range("E43").Value = "1704"
$ws.Range("F43").Value = 27578
$ws.Range("G43").Value = 689454
$ws.Range("C44").Value = "1044800187"
$ws.Range("D44").Value = "OSCAR IVAN IMITOLA HERRERA"
$ws.Range("E44").Value = "1705"
$ws.Range("F44").Value = 120000
$ws.Range("G44").Value = 3000000
$ws.Range("C45").Value = "1047397243"
$ws.Range("D45").Value = "MAIRA ALEJANDRA MARTINEZ CASTELLAR"
$ws.Range("E45").Value = "1705"
$ws.Range("F45").Value = 27578
$ws.Range("G45").Value = 689454
$ws.Range("C46").Value = "1050963899"
$ws.Range("D46").Value = "MARIA JOSE ROJAS HURTADO"
$ws.Range("E46").Value = "1705"
$ws.Range("F46").Value = 27578
$ws.Range("G46").Value = 689454
$ws.Range("C47").Value = "1044800187"
$ws.Range("D47").Value = "OSCAR IVAN IMITOLA HERRERA"
$ws.Range("E47").Value = "1706"
$ws.Range("F47").Value = 120000
$ws.Range("G47").Value = 3000000
$ws.Range("C48").Value = "1047397243"
$ws.Range("D48").Value = "MAIRA ALEJANDRA MARTINEZ CASTELLAR"
$ws.Range("E48").Value = "1706"
$ws.Range("F48").Value = 27578
$ws.Range("G48").Value = 689454
$ws.Range("C49").Value = "1050963899"
$ws.Range("D49").Value = "MARIA JOSE ROJAS HURTADO"
$ws.Range("E49").Value = "1706"
$ws.Range("F49").Value = 27578
$ws.Range("G49").Value = 689454
$ws.Range("C50").Value = "1044800187"
$ws.Range("D50").Value = "OSCAR IVAN IMITOLA HERRERA"
$ws.Range("E50").Value = "1707"
$ws.Range("F50").Value = 120000
$ws.Range("G50").Value = 3000000
$ws.Range("C51").Value = "1047397243"
$ws.Range("D51").Value = "MAIRA ALEJANDRA MARTINEZ CASTELLAR"
$ws.Range("E51").Value = "1707"
$ws.Range("F51").Value = 27578
$ws.Range("G51").Value = 689454
$ws.Range("C52").Value = "1050963899"
$ws.Range("D52").Value = "MARIA JOSE ROJAS HURTADO"
$ws.Range("E52").Value = "1707"
$ws.Range("F52").Value = 27578
$ws.Range("G52").Value = 689454
